$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H20").Value = 1500
$ws_ALC.Range("I20").Value = 1500
$ws_ALC.Range("K20").Value = 1500
$ws_ALC.Range("M20").Value = -1270
$ws_ALC.Range("H35").Value = 1500
$ws_ALC.Range("I35").Value = 1500
$ws_ALC.Range("K35").Value = 1500
$ws_ALC.Range("M35").Value = -1121
$ws_ALC.Range("H55").Value = 747.8333
$ws_ALC.Range("J55").Value = 855.8
$ws_ALC.Range("L55").Value = 855.8
$ws_ALC.Range("N55").Value = -1283.8
$ws_ALC.Range("H138").Value = 2496.0833
$ws_ALC.Range("J138").Value = 3039.2222
$ws_ALC.Range("L138").Value = 9117.6666
$ws_ALC.Range("N138").Value = -19397.6666
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H74").Value = 1463.2
$ws_ARM.Range("I74").Value = 1329.125
$ws_ARM.Range("K74").Value = 1329.125
$ws_ARM.Range("M74").Value = -455.125
$ws_ARM.Range("H77").Value = 1463.2
$ws_ARM.Range("I77").Value = 1329.125
$ws_ARM.Range("K77").Value = 6645.625
$ws_ARM.Range("M77").Value = -2277.625
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H7").Value = 10309.5
$ws_BSM.Range("I7").Value = 0
$ws_BSM.Range("J7").Value = 10309.5
$ws_BSM.Range("K7").Value = 0
$ws_BSM.Range("L7").Value = 10309.5
$ws_BSM.Range("M7").ClearContents()
$ws_BSM.Range("N7").Value = -10535.5
$ws_BSM.Range("H30").Value = 80000
$ws_BSM.Range("J30").Value = 80000
$ws_BSM.Range("L30").Value = 80000
$ws_BSM.Range("N30").Value = -80250
$ws_BSM.Range("H33").Value = 14673.333
$ws_BSM.Range("I33").Value = 14673.333
$ws_BSM.Range("K33").Value = 14673.333
$ws_BSM.Range("M33").Value = -14337.333
$ws_BSM.Range("H38").Value = 0
$ws_BSM.Range("I38").Value = 0
$ws_BSM.Range("K38").Value = 0
$ws_BSM.Range("M38").ClearContents()
$ws_BSM.Range("H99").Value = 2725.111
$ws_BSM.Range("I99").Value = 2900.125
$ws_BSM.Range("K99").Value = 2900.125
$ws_BSM.Range("M99").Value = -1402.125
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H4").Value = 20080000
$ws_CRP.Range("H6").Value = 20000020
$ws_CRP.Range("I6").Value = 20000020
$ws_CRP.Range("K6").Value = 20000020
$ws_CRP.Range("M6").Value = -19999907
$ws_CRP.Range("H19").Value = 4616861.5
$ws_CRP.Range("I19").Value = 5000183
$ws_CRP.Range("J19").Value = 17000
$ws_CRP.Range("K19").Value = 5000183
$ws_CRP.Range("L19").Value = 17000
$ws_CRP.Range("M19").Value = -5000013
$ws_CRP.Range("N19").Value = -17340
$ws_CRP.Range("H24").Value = 4616861.5
$ws_CRP.Range("I24").Value = 5000183
$ws_CRP.Range("J24").Value = 17000
$ws_CRP.Range("K24").Value = 5000183
$ws_CRP.Range("L24").Value = 17000
$ws_CRP.Range("M24").Value = -5000013
$ws_CRP.Range("N24").Value = -17340
$ws_CRP.Range("H31").Value = 1191.4
$ws_CRP.Range("I31").Value = 1195.5
$ws_CRP.Range("J31").Value = 1175
$ws_CRP.Range("K31").Value = 1195.5
$ws_CRP.Range("L31").Value = 1175
$ws_CRP.Range("M31").Value = -900.5
$ws_CRP.Range("N31").Value = -1765
$ws_CRP.Range("H34").Value = 1191.4
$ws_CRP.Range("I34").Value = 1195.5
$ws_CRP.Range("J34").Value = 1175
$ws_CRP.Range("K34").Value = 1195.5
$ws_CRP.Range("L34").Value = 1175
$ws_CRP.Range("M34").Value = -993.5
$ws_CRP.Range("N34").Value = -1579
$ws_CRP.Range("H35").Value = 5399.6665
$ws_CRP.Range("I35").Value = 5399.6665
$ws_CRP.Range("K35").Value = 5399.6665
$ws_CRP.Range("M35").Value = -5105.6665
$ws_CRP.Range("H38").Value = 22679.334
$ws_CRP.Range("I38").Value = 9019
$ws_CRP.Range("J38").Value = 50000
$ws_CRP.Range("K38").Value = 9019
$ws_CRP.Range("L38").Value = 50000
$ws_CRP.Range("M38").Value = -8642
$ws_CRP.Range("N38").Value = -50754
$ws_CRP.Range("H46").Value = 22679.334
$ws_CRP.Range("I46").Value = 9019
$ws_CRP.Range("J46").Value = 50000
$ws_CRP.Range("K46").Value = 9019
$ws_CRP.Range("L46").Value = 50000
$ws_CRP.Range("M46").Value = -8808
$ws_CRP.Range("N46").Value = -50422
$ws_CRP.Range("H132").Value = 1476.3334
$ws_CRP.Range("J132").Value = 1900
$ws_CRP.Range("L132").Value = 5700
$ws_CRP.Range("N132").Value = -10760
$ws_CRP.Range("H134").Value = 1176.8182
$ws_CRP.Range("I134").Value = 1193.1
$ws_CRP.Range("K134").Value = 3579.3
$ws_CRP.Range("M134").Value = -1044.3
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H4").Value = 2745.5715
$ws_CUL.Range("I4").Value = 2004.3182
$ws_CUL.Range("J4").Value = 4000
$ws_CUL.Range("K4").Value = 6012.9546
$ws_CUL.Range("L4").Value = 12000
$ws_CUL.Range("M4").Value = -5900.9546
$ws_CUL.Range("N4").Value = -12224
$ws_CUL.Range("H33").Value = 806.6667
$ws_CUL.Range("I33").Value = 996.6667
$ws_CUL.Range("J33").Value = 616.6667
$ws_CUL.Range("K33").Value = 5980.0002
$ws_CUL.Range("L33").Value = 3700.0002
$ws_CUL.Range("M33").Value = -5697.0002
$ws_CUL.Range("N33").Value = -4266.0002
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H5").Value = 2000
$ws_GSM.Range("I5").Value = 2000
$ws_GSM.Range("K5").Value = 2000
$ws_GSM.Range("M5").Value = -1888
$ws_GSM.Range("H132").Value = 953
$ws_GSM.Range("J132").Value = 1500
$ws_GSM.Range("L132").Value = 4500
$ws_GSM.Range("N132").Value = -9560
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H2").Value = 0
$ws_LTW.Range("J2").Value = 0
$ws_LTW.Range("L2").Value = 0
$ws_LTW.Range("N2").ClearContents()
$ws_LTW.Range("H22").Value = 3567.5
$ws_LTW.Range("I22").Value = 2650
$ws_LTW.Range("J22").Value = 3751
$ws_LTW.Range("K22").Value = 2650
$ws_LTW.Range("L22").Value = 3751
$ws_LTW.Range("M22").Value = -2355
$ws_LTW.Range("N22").Value = -4341
$ws_LTW.Range("H27").Value = 3567.5
$ws_LTW.Range("I27").Value = 2650
$ws_LTW.Range("J27").Value = 3751
$ws_LTW.Range("K27").Value = 2650
$ws_LTW.Range("L27").Value = 3751
$ws_LTW.Range("M27").Value = -2543
$ws_LTW.Range("N27").Value = -3965
$ws_LTW.Range("H61").Value = 2579
$ws_LTW.Range("I61").Value = 3010.8
$ws_LTW.Range("J61").Value = 1499.5
$ws_LTW.Range("K61").Value = 3010.8
$ws_LTW.Range("L61").Value = 1499.5
$ws_LTW.Range("M61").Value = -2808.8
$ws_LTW.Range("N61").Value = -1903.5
$ws_LTW.Range("H93").Value = 1629
$ws_LTW.Range("I93").Value = 1640.7273
$ws_LTW.Range("J93").Value = 1500
$ws_LTW.Range("K93").Value = 1640.7273
$ws_LTW.Range("L93").Value = 1500
$ws_LTW.Range("M93").Value = -392.7273
$ws_LTW.Range("N93").Value = -3996
$ws_LTW.Range("H113").Value = 2579
$ws_LTW.Range("I113").Value = 3010.8
$ws_LTW.Range("J113").Value = 1499.5
$ws_LTW.Range("K113").Value = 3010.8
$ws_LTW.Range("L113").Value = 1499.5
$ws_LTW.Range("M113").Value = -840.8000000000002
$ws_LTW.Range("N113").Value = -5839.5
$ws_LTW.Range("H132").Value = 22313.8
$ws_LTW.Range("I132").Value = 22313.8
$ws_LTW.Range("K132").Value = 66941.39999999999
$ws_LTW.Range("M132").Value = -64411.39999999999
$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H2").Value = 5000
$ws_WVR.Range("J2").Value = 5000
$ws_WVR.Range("L2").Value = 5000
$ws_WVR.Range("N2").Value = -5224
$ws_WVR.Range("H5").Value = 4600000
$ws_WVR.Range("J5").Value = 3250000
$ws_WVR.Range("L5").Value = 3250000
$ws_WVR.Range("N5").Value = -3250224
$ws_WVR.Range("H132").Value = 1952.1428
$ws_WVR.Range("J132").Value = 2498.25
$ws_WVR.Range("L132").Value = 7494.75
$ws_WVR.Range("N132").Value = -12554.75
